$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet stripes its data rows by absolute row position (even rows use
# one style family, odd rows use another) rather than by row content.
# Inserting a row shifts the old rows' content AND their old formatting
# down together, which breaks the stripe pattern below the insertion
# point. So: stash the two alternating formats on a scratch sheet (a plain
# row-insert on the data sheet would otherwise drag any same-sheet "parked"
# template cells down with it), insert the new row, then re-apply the
# correct alternating format to every data row from the new row down.
$scratch = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$ws.Range("A2:G2").Copy()
$scratch.Range("A1:G1").PasteSpecial(-4122)  # xlPasteFormats -> "even" row template
$ws.Range("A3:G3").Copy()
$scratch.Range("A2:G2").PasteSpecial(-4122)  # xlPasteFormats -> "odd" row template
$ws.Application.CutCopyMode = 0

# Insert a new row above row 2. Existing rows 2-7 (PED-B2-1..PED-B2-6) shift
# down to rows 3-8, unchanged in content.
$ws.Rows.Item(2).Insert()

# Re-stripe every data row (2-8) by its new position parity.
for ($r = 2; $r -le 8; $r++) {
  if ($r % 2 -eq 0) {
    $scratch.Range("A1:G1").Copy()
  } else {
    $scratch.Range("A2:G2").Copy()
  }
  $ws.Range("A" + $r + ":G" + $r).PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = 0

$null = $scratch.Delete()

# Fill in the new quiz session for PED-B2-1 (a second session, one day
# later). Force the date/time columns to text before writing so the
# date/time-shaped strings are kept as literal text -- matching how the
# rest of the sheet stores these columns -- instead of being auto-converted
# to Excel date/time serials.
$dateFmt = $ws.Range("E2").NumberFormat
$timeFmt = $ws.Range("F2").NumberFormat

$ws.Range("A2").Value2 = "Year 4"
$ws.Range("B2").Value2 = "PED-B2-1"
$ws.Range("C2").Value2 = "pediatrics"
$ws.Range("D2").Value2 = "2"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "07/09/2025"
$ws.Range("E2").NumberFormat = $dateFmt

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value2 = "08:00:00"
$ws.Range("F2").NumberFormat = $timeFmt

$ws.Range("G2").Value2 = 120
